$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated weights (column B), optimal portfolio (column C) and
# optimal portfolio with view (column D) values.
$ws.Range("B2").Value = 0.1
$ws.Range("C2").Value = 0.1272849253007465
$ws.Range("D2").Value = 0

$ws.Range("B3").Value = 0.05
$ws.Range("C3").Value = 0.1272908224383453
$ws.Range("D3").Value = 0.1163737759267746

$ws.Range("B4").Value = 0.1
$ws.Range("C4").Value = 0.1796008618156742
$ws.Range("D4").Value = 0.2176989893803725

$ws.Range("B5").Value = 0.1
$ws.Range("C5").Value = 0.1722192567577291
$ws.Range("D5").Value = 0.1949831861608918

$ws.Range("B6").Value = 0.15
$ws.Range("C6").Value = 0.1389108338040999
$ws.Range("D6").Value = 0.07649963773174576

$ws.Range("B7").Value = 0.2
$ws.Range("C7").Value = 0.1272876650722808
$ws.Range("D7").Value = 0.223703413766908

$ws.Range("B8").Value = 0.3
$ws.Range("C8").Value = 0.1274056348111242
$ws.Range("D8").Value = 0.1707409970333074
